# Update the "想去人数" (want-to-go count) figures for a handful of events
# on the 展览 (Exhibition) and 全部类型 (All Types) sheets, reflecting the
# latest scrape of the gh-pages generated output.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # 蜀山·银泰百货高新店-2024漫趣地带嘉年华（免费）
    $ws.Range("F2").Value = 295

    # 合肥·第十五届次元之门动漫游戏博览会
    $ws.Range("F4").Value = 7928

    # 合肥·首届AT次元时代动漫游戏嘉年华
    $ws.Range("F5").Value = 5783
}

# 合肥·第九届环形宇宙动漫游戏嘉年华 is at row 11 on 展览 and row 14 on 全部类型
$wb.Worksheets.Item("展览").Range("F11").Value = 330
$wb.Worksheets.Item("全部类型").Range("F14").Value = 330
